$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.044.96'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.828.87'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.84'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6224'
$ws.Range('E6').Value = '  -6.19%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.48'
$ws.Range('E8').Value = '  +6.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07388'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2921'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.68'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07601'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('D13').Value = '1.829.64'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6626'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.11'
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009082'
$ws.Range('E17').Value = '  +8.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.010'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = '29.046.69'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = '2.079.44'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '225.19'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.35'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.166'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.29'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.414'
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1355'
$ws.Range('E28').Value = '  -3.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.80'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.497'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.027'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05239'
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.837'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7348'
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.150'
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '1.280.09'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01782'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.324'
$ws.Range('E42').Value = '  +6.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8973'
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.56'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').Value = '1.977.61'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.66'
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000120'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.705'
$ws.Range('E50').Value = '  -3.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3962'
$ws.Range('E51').Value = '  -1.34%  '
